# Apply corrected Relevance Marker values (Appenzeller-Herzog 2019 - van Dis 2020)
# for the time-to-discovery simulation results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C, rows 2-98 (td_sim_1)
$cValues = @(
    405, 22, 291, 91, 27, 260, 146, 124, 18, 103,
    119, 71, 73, 88, 210, 90, 50, 30, 356, 5,
    35, 12, 60, 108, 38, 40, 87, 10, 26, 278,
    332, 150, 63, 78, 25, 1102, 502, 472, 218, 338,
    224, 17, 4, 64, 354, 19, 300, 232, 155, 45,
    440, 9, 118, 16, 325, 28, 312, 117, 72, 15,
    485, 256, 122, 46, 104, 190, 3, 39, 70, 62,
    198, 77, 114, 139, 471, 23, 110, 93, 116, 468,
    89, 448, 469, 8, 445, 151, 259, 86, 75, 457,
    61, 185, 11, 566, 138, 82, 168.0729166666667
)

# New values for column D, rows 2-97 (average_simulation_TD)
$dValues = @(
    407, 20.5, 299.5, 78, 24, 256.5, 153, 142.5, 15, 106,
    112.5, 69, 101, 80, 193.5, 82, 54.5, 27, 355.5, 3,
    30, 11, 56, 112, 36, 34, 79, 10, 24.5, 235.5,
    332.5, 143, 47, 80.5, 32.5, 1100, 498, 464, 224.5, 331.5,
    232, 17.5, 3, 62.5, 364, 16, 285, 260.5, 154, 42,
    446.5, 7.5, 111, 15.5, 328.5, 30.5, 293, 109.5, 70, 12,
    502.5, 250.5, 120, 38, 107.5, 199.5, 4, 37, 60.5, 49,
    194, 71.5, 106.5, 143, 458, 21.5, 103.5, 87, 108.5, 456,
    81, 452.5, 455.5, 6, 452.5, 151.5, 255.5, 75, 68.5, 474,
    53, 207.5, 9.5, 560.5, 127.5, 75.5
)

$cStartRow = 2
for ($i = 0; $i -lt $cValues.Length; $i++) {
    $ws.Cells.Item($cStartRow + $i, 3).Value2 = $cValues[$i]
}

$dStartRow = 2
for ($i = 0; $i -lt $dValues.Length; $i++) {
    $ws.Cells.Item($dStartRow + $i, 4).Value2 = $dValues[$i]
}

Write-Output "Updated $($cValues.Length) C-column cells and $($dValues.Length) D-column cells."
